$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1090.3572
$ws.Range("J19").Value = 1579.5333
$ws.Range("L19").Value = 1579.5333
$ws.Range("N19").Value = -1929.5333
$ws.Range("H33").Value = 4670.8335
$ws.Range("I33").Value = 5223.625
$ws.Range("K33").Value = 5223.625
$ws.Range("M33").Value = -4994.625
$ws.Range("H53").Value = 4964.3477
$ws.Range("J53").Value = 9305.916999999999
$ws.Range("L53").Value = 9305.916999999999
$ws.Range("N53").Value = -10579.917
$ws.Range("H62").Value = 3207.4614
$ws.Range("I62").Value = 3207.4614
$ws.Range("K62").Value = 3207.4614
$ws.Range("M62").Value = -2583.4614
$ws.Range("H65").Value = 3207.4614
$ws.Range("I65").Value = 3207.4614
$ws.Range("K65").Value = 16037.307
$ws.Range("M65").Value = -12917.307
$ws.Range("H88").Value = 3860.3
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3860.3
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3860.3
$ws.Range("M88").Value = ""
$ws.Range("N88").Value = -4672.3
$ws.Range("H91").Value = 3860.3
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3860.3
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3860.3
$ws.Range("M91").Value = ""
$ws.Range("N91").Value = -6668.3
$ws.Range("H127").Value = 1849.091
$ws.Range("I127").Value = 1442.5
$ws.Range("K127").Value = 4327.5
$ws.Range("M127").Value = 632.5
$ws.Range("H132").Value = 40003280
$ws.Range("I132").Value = 47622620
$ws.Range("K132").Value = 142867860
$ws.Range("M132").Value = -142865330
$ws.Range("H137").Value = 150505.67
$ws.Range("I137").Value = 223971.88
$ws.Range("K137").Value = 671915.64
$ws.Range("M137").Value = -669365.64
$ws.Range("H138").Value = 2255.4285
$ws.Range("J138").Value = 4020.8965
$ws.Range("L138").Value = 12062.6895
$ws.Range("N138").Value = -22342.6895
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5584.1177
$ws.Range("I32").Value = 3268.513
$ws.Range("K32").Value = 3268.513
$ws.Range("M32").Value = -2981.513
$ws.Range("H76").Value = 73000.336
$ws.Range("J76").Value = 73000.336
$ws.Range("L76").Value = 73000.336
$ws.Range("N76").Value = -73676.336
$ws.Range("H79").Value = 73000.336
$ws.Range("J79").Value = 73000.336
$ws.Range("L79").Value = 73000.336
$ws.Range("N79").Value = -75340.336
$ws.Range("H122").Value = 523099.5
$ws.Range("I122").Value = 1948.258
$ws.Range("K122").Value = 5844.774
$ws.Range("M122").Value = -3394.774
$ws.Range("H132").Value = 2833.2307
$ws.Range("I132").Value = 2351.25
$ws.Range("K132").Value = 7053.75
$ws.Range("M132").Value = -4523.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2500.6
$ws.Range("I20").Value = 1667.3334
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 1667.3334
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -1420.3334
$ws.Range("N20").Value = -10494
$ws.Range("H94").Value = 3796660.2
$ws.Range("I94").Value = 6061835
$ws.Range("J94").Value = 21368.889
$ws.Range("K94").Value = 6061835
$ws.Range("L94").Value = 21368.889
$ws.Range("M94").Value = -6061384
$ws.Range("N94").Value = -22270.889
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14063.269
$ws.Range("I31").Value = 2257.6155
$ws.Range("K31").Value = 2257.6155
$ws.Range("M31").Value = -1962.6155
$ws.Range("H34").Value = 14063.269
$ws.Range("I34").Value = 2257.6155
$ws.Range("K34").Value = 2257.6155
$ws.Range("M34").Value = -2055.6155
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""
$ws.Range("H132").Value = 39019
$ws.Range("I132").Value = 2169.8572
$ws.Range("K132").Value = 6509.571599999999
$ws.Range("M132").Value = -3979.571599999999
$ws.Range("H134").Value = 2502.7932
$ws.Range("I134").Value = 1743.65
$ws.Range("K134").Value = 5230.950000000001
$ws.Range("M134").Value = -2695.950000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 349.16666
$ws.Range("I2").Value = 54.5
$ws.Range("J2").Value = 456.31818
$ws.Range("K2").Value = 327
$ws.Range("L2").Value = 2737.90908
$ws.Range("M2").Value = -214
$ws.Range("N2").Value = -2963.90908
$ws.Range("H56").Value = 125005600
$ws.Range("I56").Value = 125005600
$ws.Range("K56").Value = 125005600
$ws.Range("M56").Value = -125005070
$ws.Range("H98").Value = 503.57144
$ws.Range("I98").Value = 446.5
$ws.Range("K98").Value = 1339.5
$ws.Range("M98").Value = 158.5
$ws.Range("H130").Value = 3500
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").Value = ""
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6252.25
$ws.Range("I12").Value = 5003
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 5003
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -4863
$ws.Range("N12").Value = -10280
$ws.Range("H122").Value = 225495.53
$ws.Range("I122").Value = 298896.12
$ws.Range("K122").Value = 896688.36
$ws.Range("M122").Value = -894238.36
$ws.Range("H126").Value = 7578799.5
$ws.Range("J126").Value = 20836286
$ws.Range("L126").Value = 62508858
$ws.Range("N126").Value = -62513798
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5447.4
$ws.Range("I7").Value = 2175.6667
$ws.Range("K7").Value = 2175.6667
$ws.Range("M7").Value = -2063.6667
$ws.Range("H61").Value = 10105568
$ws.Range("I61").Value = 11115941
$ws.Range("J61").Value = 1845
$ws.Range("K61").Value = 11115941
$ws.Range("L61").Value = 1845
$ws.Range("M61").Value = -11115739
$ws.Range("N61").Value = -2249
$ws.Range("H113").Value = 10105568
$ws.Range("I113").Value = 11115941
$ws.Range("J113").Value = 1845
$ws.Range("K113").Value = 11115941
$ws.Range("L113").Value = 1845
$ws.Range("M113").Value = -11113771
$ws.Range("N113").Value = -6185
$ws.Range("H126").Value = 5447.4
$ws.Range("I126").Value = 2175.6667
$ws.Range("K126").Value = 6527.000100000001
$ws.Range("M126").Value = -4057.000100000001
$ws.Range("H132").Value = 8132.8423
$ws.Range("I132").Value = 8303.808000000001
$ws.Range("J132").Value = 7762.4165
$ws.Range("K132").Value = 24911.424
$ws.Range("L132").Value = 23287.2495
$ws.Range("M132").Value = -22381.424
$ws.Range("N132").Value = -28347.2495
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8109.0938
$ws.Range("I62").Value = 4400.8
$ws.Range("J62").Value = 8795.814
$ws.Range("K62").Value = 4400.8
$ws.Range("L62").Value = 8795.814
$ws.Range("M62").Value = -3776.8
$ws.Range("N62").Value = -10043.814
$ws.Range("H65").Value = 8109.0938
$ws.Range("I65").Value = 4400.8
$ws.Range("J65").Value = 8795.814
$ws.Range("K65").Value = 22004
$ws.Range("L65").Value = 43979.07
$ws.Range("M65").Value = -18884
$ws.Range("N65").Value = -50219.07
$ws.Range("H122").Value = 2897.6
$ws.Range("I122").Value = 1798.4
$ws.Range("J122").Value = 3996.8
$ws.Range("K122").Value = 5395.200000000001
$ws.Range("L122").Value = 11990.4
$ws.Range("M122").Value = -2945.200000000001
$ws.Range("N122").Value = -16890.4
$ws.Range("H126").Value = 3666.8
$ws.Range("I126").Value = 3182.0908
$ws.Range("J126").Value = 4999.75
$ws.Range("K126").Value = 9546.2724
$ws.Range("L126").Value = 14999.25
$ws.Range("M126").Value = -7076.2724
$ws.Range("N126").Value = -19939.25
$ws.Range("H132").Value = 47668450
$ws.Range("I132").Value = 50007424
$ws.Range("K132").Value = 150022272
$ws.Range("M132").Value = -150019742

Write-Output "Applied all updates"